# Update "Förändrad" (column C) dates from 45184 -> 45186 on every data row,
# and append the friendly display name (the value of column A on that row)
# as the second HYPERLINK() argument for any S/T/V/W/X/Y cell that still
# has the old single-argument form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$firstRow = $ur.Row
$lastRow = $firstRow + $ur.Rows.Count - 1

# Data starts on row 2 (row 1 is the header).
$startRow = 2
if ($firstRow -gt $startRow) { $startRow = $firstRow }

for ($r = $startRow; $r -le $lastRow; $r++) {

    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    $label = $ws.Cells.Item($r, 1).Value2

    foreach ($col in 19, 20, 22, 23, 24, 25) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f.Contains("HYPERLINK(") -and -not $f.Contains(",")) {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $label + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
